$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7602039999999999
$ws.Range("H2").Value = 2.280612
$ws.Range("I2").Value = 0.9081302434927829
$ws.Range("J2").Value = 0.9140443705567521
$ws.Range("M2").Value = 8.021311666666668
$ws.Range("N2").Value = 24.063935
$ws.Range("O2").Value = 0.2318437811880539
$ws.Range("P2").Value = 0.2575223477274881
$ws.Range("Q2").Value = 6.097833214246666
$ws.Range("R2").Value = 54.88049892821999
$ws.Range("S2").Value = 0.2105443494625948
$ws.Range("T2").Value = 0.2353868522328689
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7602039999999999
$ws.Range("H3").Value = 2.280612
$ws.Range("I3").Value = 0.9081302434927829
$ws.Range("J3").Value = 0.9140443705567521
$ws.Range("O3").Value = 0.3382003894878138
$ws.Range("P3").Value = 0.3756588072233373
$ws.Range("Q3").Value = 8.895168796514664
$ws.Range("R3").Value = 80.05651916863198
$ws.Range("S3").Value = 0.3071300020549224
$ws.Range("T3").Value = 0.3433688179925556
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7602039999999999
$ws.Range("H4").Value = 2.280612
$ws.Range("I4").Value = 0.9081302434927829
$ws.Range("J4").Value = 0.9140443705567521
$ws.Range("M4").Value = 2.475956
$ws.Range("N4").Value = 7.427868
$ws.Range("O4").Value = 0.07156373233578578
$ws.Range("P4").Value = 0.07948999222154987
$ws.Range("Q4").Value = 1.882231655024
$ws.Range("R4").Value = 16.940084895216
$ws.Range("S4").Value = 0.06498918967134948
$ws.Range("T4").Value = 0.07265737990570767
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7602039999999999
$ws.Range("H5").Value = 2.280612
$ws.Range("I5").Value = 0.9081302434927829
$ws.Range("J5").Value = 0.9140443705567521
$ws.Range("M5").Value = 10.349683
$ws.Range("N5").Value = 20.699366
$ws.Range("O5").Value = 0.2991418038011306
$ws.Range("P5").Value = 0.2215161123395049
$ws.Range("Q5").Value = 7.867870415331998
$ws.Range("R5").Value = 47.20722249199198
$ws.Range("S5").Value = 0.2716597191247911
$ws.Range("T5").Value = 0.2024755554715416
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7602039999999999
$ws.Range("H6").Value = 2.280612
$ws.Range("I6").Value = 0.9081302434927829
$ws.Range("J6").Value = 0.9140443705567521
$ws.Range("M6").Value = 2.049936666666667
$ws.Range("N6").Value = 6.14981
$ws.Range("O6").Value = 0.05925029318721588
$ws.Range("P6").Value = 0.06581274048811983
$ws.Range("Q6").Value = 1.558370053746667
$ws.Range("R6").Value = 14.02533048372
$ws.Range("S6").Value = 0.05380698317912513
$ws.Range("T6").Value = 0.06015576495407837
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.016249
$ws.Range("H7").Value = 0.032498
$ws.Range("I7").Value = 0.01941085330584189
$ws.Range("J7").Value = 0.01302484331151171
$ws.Range("M7").Value = 8.021311666666668
$ws.Range("N7").Value = 24.063935
$ws.Range("O7").Value = 0.2318437811880539
$ws.Range("P7").Value = 0.2575223477274881
$ws.Range("Q7").Value = 0.1303382932716667
$ws.Range("R7").Value = 0.78202975963
$ws.Range("S7").Value = 0.00450028562651302
$ws.Range("T7").Value = 0.003354188228363165
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.016249
$ws.Range("H8").Value = 0.032498
$ws.Range("I8").Value = 0.01941085330584189
$ws.Range("J8").Value = 0.01302484331151171
$ws.Range("O8").Value = 0.3382003894878138
$ws.Range("P8").Value = 0.3756588072233373
$ws.Range("Q8").Value = 0.1901300148046667
$ws.Range("R8").Value = 1.140780088828
$ws.Range("S8").Value = 0.006564758148326547
$ws.Range("T8").Value = 0.00489289710267335
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.016249
$ws.Range("H9").Value = 0.032498
$ws.Range("I9").Value = 0.01941085330584189
$ws.Range("J9").Value = 0.01302484331151171
$ws.Range("M9").Value = 2.475956
$ws.Range("N9").Value = 7.427868
$ws.Range("O9").Value = 0.07156373233578578
$ws.Range("P9").Value = 0.07948999222154987
$ws.Range("Q9").Value = 0.040231809044
$ws.Range("R9").Value = 0.241390854264
$ws.Range("S9").Value = 0.001389113110388472
$ws.Range("T9").Value = 0.001035344693518972
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.016249
$ws.Range("H10").Value = 0.032498
$ws.Range("I10").Value = 0.01941085330584189
$ws.Range("J10").Value = 0.01302484331151171
$ws.Range("M10").Value = 10.349683
$ws.Range("N10").Value = 20.699366
$ws.Range("O10").Value = 0.2991418038011306
$ws.Range("P10").Value = 0.2215161123395049
$ws.Range("Q10").Value = 0.168171999067
$ws.Range("R10").Value = 0.6726879962679999
$ws.Range("S10").Value = 0.005806597671228684
$ws.Range("T10").Value = 0.002885212654197276
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.016249
$ws.Range("H11").Value = 0.032498
$ws.Range("I11").Value = 0.01941085330584189
$ws.Range("J11").Value = 0.01302484331151171
$ws.Range("M11").Value = 2.049936666666667
$ws.Range("N11").Value = 6.14981
$ws.Range("O11").Value = 0.05925029318721588
$ws.Range("P11").Value = 0.06581274048811983
$ws.Range("Q11").Value = 0.03330942089666667
$ws.Range("R11").Value = 0.19985652538
$ws.Range("S11").Value = 0.001150098749385171
$ws.Range("T11").Value = 0.0008572006327589432
$ws.Range("G12").Value = 0.06065600000000001
$ws.Range("H12").Value = 0.181968
$ws.Range("I12").Value = 0.07245890320137523
$ws.Range("J12").Value = 0.07293078613173619
$ws.Range("M12").Value = 8.021311666666668
$ws.Range("N12").Value = 24.063935
$ws.Range("O12").Value = 0.2318437811880539
$ws.Range("P12").Value = 0.2575223477274881
$ws.Range("Q12").Value = 0.4865406804533334
$ws.Range("R12").Value = 4.378866124080001
$ws.Range("S12").Value = 0.01679914609894602
$ws.Range("T12").Value = 0.01878130726625604
$ws.Range("G13").Value = 0.06065600000000001
$ws.Range("H13").Value = 0.181968
$ws.Range("I13").Value = 0.07245890320137523
$ws.Range("J13").Value = 0.07293078613173619
$ws.Range("O13").Value = 0.3382003894878138
$ws.Range("P13").Value = 0.3756588072233373
$ws.Range("Q13").Value = 0.7097375948053334
$ws.Range("R13").Value = 6.387638353248001
$ws.Range("S13").Value = 0.0245056292845649
$ws.Range("T13").Value = 0.02739709212810832
$ws.Range("G14").Value = 0.06065600000000001
$ws.Range("H14").Value = 0.181968
$ws.Range("I14").Value = 0.07245890320137523
$ws.Range("J14").Value = 0.07293078613173619
$ws.Range("M14").Value = 2.475956
$ws.Range("N14").Value = 7.427868
$ws.Range("O14").Value = 0.07156373233578578
$ws.Range("P14").Value = 0.07948999222154987
$ws.Range("Q14").Value = 0.150181587136
$ws.Range("R14").Value = 1.351634284224
$ws.Range("S14").Value = 0.005185429554047828
$ws.Range("T14").Value = 0.005797267622323227
$ws.Range("G15").Value = 0.06065600000000001
$ws.Range("H15").Value = 0.181968
$ws.Range("I15").Value = 0.07245890320137523
$ws.Range("J15").Value = 0.07293078613173619
$ws.Range("M15").Value = 10.349683
$ws.Range("N15").Value = 20.699366
$ws.Range("O15").Value = 0.2991418038011306
$ws.Range("P15").Value = 0.2215161123395049
$ws.Range("Q15").Value = 0.627770372048
$ws.Range("R15").Value = 3.766622232288
$ws.Range("S15").Value = 0.0216754870051109
$ws.Range("T15").Value = 0.01615534421376608
$ws.Range("G16").Value = 0.06065600000000001
$ws.Range("H16").Value = 0.181968
$ws.Range("I16").Value = 0.07245890320137523
$ws.Range("J16").Value = 0.07293078613173619
$ws.Range("M16").Value = 2.049936666666667
$ws.Range("N16").Value = 6.14981
$ws.Range("O16").Value = 0.05925029318721588
$ws.Range("P16").Value = 0.06581274048811983
$ws.Range("Q16").Value = 0.1243409584533334
$ws.Range("R16").Value = 1.11906862608
$ws.Range("S16").Value = 0.004293211258705578
$ws.Range("T16").Value = 0.004799774901282522
